$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the billing line (row 3) with the new competencia figures ---
$ws.Range("D3").Value = 1000.0
$ws.Range("F3").Value = 993.85
$ws.Range("G3").Value = 0.0
$ws.Range("H3").Value = 993.85
$ws.Range("I3").Value = "17/09/2025"
$ws.Range("J3").Value = 0.0
$ws.Range("K3").Value = "RECURSAR"
$ws.Range("L3").Value = "18/09/2025"

# --- TOTAL row (row 4) now sums every money column, not just GLOSA ---
$ws.Range("D4").Formula = "=SUM(D3:D3)"
$ws.Range("E4").Formula = "=SUM(E3:E3)"
$ws.Range("F4").Formula = "=SUM(F3:F3)"
$ws.Range("G4").Formula = "=SUM(G3:G3)"
$ws.Range("H4").Formula = "=SUM(H3:H3)"
$ws.Range("J4").Formula = "=SUM(J3:J3)"

# --- Shrink the bold data/total-row fonts from 12pt to 11pt ---
$ws.Range("A3:L4").Font.Size = 11

# --- Column widths settle to their new best-fit values after the font/content change ---
$ws.Columns("A").ColumnWidth = 16.666666666666668
$ws.Columns("B").ColumnWidth = 14.333333333333334
$ws.Columns("G").ColumnWidth = 11.666666666666666
$ws.Columns("H").ColumnWidth = 13.166666666666666
$ws.Columns("I").ColumnWidth = 12.0
$ws.Columns("J").ColumnWidth = 13.166666666666666
$ws.Columns("L").ColumnWidth = 14.333333333333334

# --- Print setup: landscape, fit to 1 page wide, no longer constrained to 1 page tall ---
$ws.PageSetup.Orientation = 2
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = $false
